$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCF")
Write-Host $ws.Name
